# Update the "Metadata" sheet (StructureDefinition version bump + publisher/jurisdiction
# info) and the "Elements" sheet (root Extension's Short/Definition text) to match the
# new IG build (version 6.0.0, dated 2022-01-21).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# -- Metadata sheet --------------------------------------------------------

# Version bump: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date of this IG build
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher was blank before; now populated
$meta.Range("B9").Value = "Alvearie Team"

# The old row 10/11 pair both held "Contact" / "No display for ContactDetail".
# Row 11 (the duplicate) is removed entirely, and row 10 is repurposed to show
# the Jurisdiction property instead of the duplicated Contact property.
$meta.Rows.Item(11).Delete()

$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# -- Elements sheet ---------------------------------------------------------

# Root "Extension" element's Short/Definition now mirror the StructureDefinition's
# own Title/Description instead of the generic boilerplate text.
$elements.Range("K2").Value = "Employee Family Size"
$elements.Range("L2").Value = "Total family size of the employee, whether or not dependents are insured or participating in any programs"
